$d = $word.ActiveDocument

# 1. Shorten the dots on the "Олгосон нярав" line.
$d.Content.Find.Execute(
    "Олгосон нярав: .........................................................",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Олгосон нярав: ...........................................", 2)

# 2. Remove the old _GoBack bookmark (currently sitting in the
#    "Шаардах бичсэн" paragraph) -- it will be re-created at the very end
#    of the document, as in the target.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 3. Append the new paragraphs (3 blank ones, "Өбы", the "Шардлагатай
#    өөрчлөлт оруулсан" line, and a final blank paragraph that now owns
#    the _GoBack bookmark) at the very end of the document body.
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$xml  = "<w:p $ns><w:pPr><w:rPr><w:lang w:val=`"mn-MN`"/></w:rPr></w:pPr></w:p>"
$xml += "<w:p $ns><w:pPr><w:rPr><w:lang w:val=`"mn-MN`"/></w:rPr></w:pPr></w:p>"
$xml += "<w:p $ns><w:pPr><w:rPr><w:lang w:val=`"mn-MN`"/></w:rPr></w:pPr></w:p>"
$xml += "<w:p $ns><w:r><w:rPr><w:lang w:val=`"mn-MN`"/></w:rPr><w:t>Өбы</w:t></w:r></w:p>"
$xml += "<w:p $ns><w:pPr><w:rPr><w:lang w:val=`"mn-MN`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"mn-MN`"/></w:rPr><w:t>Шардлагатай өөрчлөлт оруулсан</w:t></w:r></w:p>"
$xml += "<w:p $ns><w:pPr><w:rPr><w:lang w:val=`"mn-MN`"/></w:rPr></w:pPr><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>"

$r = $d.Content
$r.Collapse(0)
$r.InsertXML($xml)
